$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings that must stay as TEXT
# (matching the original inline-string cell type), so we force a text number
# format before assigning, to prevent Excel auto-converting them to numbers.
$textForceCells = @(
    "D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.473.05"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "1.672.39"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "313.62"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "0.3902"
$ws.Range("E7").Value = "  -3.04%  "
$ws.Range("E8").Value = "  -3.03%  "
$ws.Range("D9").Value = "1.006"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "51.59"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("D11").Value = "1.396"
$ws.Range("E11").Value = "  -5.30%  "
$ws.Range("D12").Value = "0.08604"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "24.98"
$ws.Range("E13").Value = "  -3.92%  "
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").Value = "0.00001309"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "7.706"
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("D17").Value = "1.667.21"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "93.41"
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("D19").Value = "0.07053"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").Value = "20.40"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").Value = "7.040"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "1.007"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "13.95"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").Value = "24.474.97"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").Value = "2.372"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("D26").Value = "23.28"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").Value = "2.744"
$ws.Range("E27").Value = "  -4.88%  "
$ws.Range("D28").Value = "162.08"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").Value = "5.799"
$ws.Range("E29").Value = "  -9.57%  "
$ws.Range("D30").Value = "146.73"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "8.221"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "2.542"
$ws.Range("E32").Value = "  +13.53%  "
$ws.Range("D33").Value = "1.852.30"
$ws.Range("E33").Value = "  -5.22%  "
$ws.Range("D34").Value = "0.08331"
$ws.Range("E34").Value = "  -6.59%  "
$ws.Range("D35").Value = "0.03015"
$ws.Range("E35").Value = "  -5.61%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "6.930"
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2793"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").Value = "0.9656"
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("D39").Value = "0.09478"
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("D40").Value = "1.514"
$ws.Range("E40").Value = "  +3.15%  "
$ws.Range("D41").Value = "10.26"
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("D42").Value = "0.7846"
$ws.Range("E42").Value = "  -6.65%  "
$ws.Range("D43").Value = "13.47"
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("D44").Value = "16.36"
$ws.Range("E44").Value = "  -6.15%  "
$ws.Range("D45").Value = "0.7069"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("D46").Value = "2.543"
$ws.Range("E46").Value = "  -5.85%  "
$ws.Range("D47").Value = "4.174"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.08581"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("D50").Value = "1.313"
$ws.Range("E50").Value = "  -5.14%  "
$ws.Range("D51").Value = "136.98"
$ws.Range("E51").Value = "  -3.49%  "
